$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to hold literal TEXT (not auto-converted to a number)
# by temporarily marking the cell as Text-formatted, assigning the value,
# then resetting the style back to "Normal" so no stray formatting is left behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "91.994.61"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.081.59"
$ws.Range("E3").Value = "  -2.57%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
Set-TextValue $ws.Range("D5") "236.40"
$ws.Range("E5").Value = "  -1.33%  "

# Row 6
Set-TextValue $ws.Range("D6") "608.19"
$ws.Range("E6").Value = "  -2.06%  "

# Row 7
Set-TextValue $ws.Range("D7") "1.08"
$ws.Range("E7").Value = "  -4.62%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.387"
$ws.Range("E8").Value = "  +2.91%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.999"
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
Set-TextValue $ws.Range("D10") "3.075.34"
$ws.Range("E10").Value = "  -2.76%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.730"
$ws.Range("E11").Value = "  -2.13%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.199"
$ws.Range("E12").Value = "  -2.10%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000246"
$ws.Range("E13").Value = "  -0.72%  "

# Row 14
Set-TextValue $ws.Range("D14") "92.120.10"
$ws.Range("E14").Value = "  +0.96%  "

# Row 15
Set-TextValue $ws.Range("D15") "33.72"
$ws.Range("E15").Value = "  -5.16%  "

# Row 16
Set-TextValue $ws.Range("D16") "5.39"
$ws.Range("E16").Value = "  -3.38%  "

# Row 17
Set-TextValue $ws.Range("D17") "3.668.65"
$ws.Range("E17").Value = "  -2.24%  "

# Row 18
Set-TextValue $ws.Range("D18") "3.104.62"
$ws.Range("E18").Value = "  -2.62%  "

# Row 19
Set-TextValue $ws.Range("D19") "3.73"
$ws.Range("E19").Value = "  -0.52%  "

# Row 20
Set-TextValue $ws.Range("D20") "14.52"
$ws.Range("E20").Value = "  -4.62%  "

# Row 21
Set-TextValue $ws.Range("D21") "5.68"
$ws.Range("E21").Value = "  -4.13%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D22") "9.25"
$ws.Range("E22").Value = "  +0.41%  "

# Row 23
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D23") "440.77"
$ws.Range("E23").Value = "  -3.79%  "

# Row 24
Set-TextValue $ws.Range("D24") "0.0000192"
$ws.Range("E24").Value = "  -5.48%  "

# Row 25
Set-TextValue $ws.Range("D25") "5.63"
$ws.Range("E25").Value = "  -6.91%  "

# Row 26
Set-TextValue $ws.Range("D26") "85.46"
$ws.Range("E26").Value = "  -3.84%  "

# Row 27
Set-TextValue $ws.Range("D27") "11.51"
$ws.Range("E27").Value = "  -4.59%  "

# Row 28
Set-TextValue $ws.Range("D28") "3.253.20"
$ws.Range("E28").Value = "  -1.97%  "

# Row 29
Set-TextValue $ws.Range("D29") "0.995"
$ws.Range("E29").Value = "  -0.22%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.132"
$ws.Range("E30").Value = "  +4.25%  "

# Row 31
Set-TextValue $ws.Range("D31") "0.226"
$ws.Range("E31").Value = "  -2.57%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.167"
$ws.Range("E32").Value = "  -2.94%  "

# Row 33
Set-TextValue $ws.Range("D33") "8.99"
$ws.Range("E33").Value = "  -4.36%  "

# Row 34
$ws.Range("E34").Value = "  +6.58%  "

# Row 35
Set-TextValue $ws.Range("D35") "7.80"
$ws.Range("E35").Value = "  +1.01%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.157"
$ws.Range("E36").Value = "  -8.02%  "

# Row 37
Set-TextValue $ws.Range("D37") "25.71"
$ws.Range("E37").Value = "  -3.66%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.89"
$ws.Range("E38").Value = "  +1.07%  "

# Row 39
$ws.Range("E39").Value = "  -3.18%  "

# Row 40
Set-TextValue $ws.Range("D40") "475.68"
$ws.Range("E40").Value = "  -7.24%  "

# Row 41
Set-TextValue $ws.Range("D41") "23.91"
$ws.Range("E41").Value = "  +7.70%  "

# Row 42
Set-TextValue $ws.Range("D42") "1.27"
$ws.Range("E42").Value = "  -6.53%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.426"
$ws.Range("E43").Value = "  -5.66%  "

# Row 44
Set-TextValue $ws.Range("D44") "3.24"
$ws.Range("E44").Value = "  -6.30%  "

# Row 45
$ws.Range("E45").Value = "  +0.08%  "

# Row 46
Set-TextValue $ws.Range("D46") "162.51"
$ws.Range("E46").Value = "  +2.27%  "

# Row 47
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D47") "0.677"
$ws.Range("E47").Value = "  -4.99%  "

# Row 48
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D48") "1.85"
$ws.Range("E48").Value = "  -5.11%  "

# Row 49
Set-TextValue $ws.Range("D49") "1.36"
$ws.Range("E49").Value = "  -1.58%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.0326"
$ws.Range("E50").Value = "  +1.30%  "

# Row 51
Set-TextValue $ws.Range("D51") "43.91"
$ws.Range("E51").Value = "  -0.41%  "
